$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Sertad1"
$ws.Cells.Item(2, 3).Value = "Ar"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 22.02847633333333
$ws.Cells.Item(2, 8).Value = 66.085429
$ws.Cells.Item(2, 9).Value = 0.3378777065804683
$ws.Cells.Item(2, 10).Value = 0.3378777065804683
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.6732596666666667
$ws.Cells.Item(2, 14).Value = 2.019779
$ws.Cells.Item(2, 15).Value = 0.05029912702290298
$ws.Cells.Item(2, 16).Value = 0.05029912702290298
$ws.Cells.Item(2, 17).Value = 14.83088463335456
$ws.Cells.Item(2, 18).Value = 133.477961700191
$ws.Cells.Item(2, 19).Value = 0.01699495368149812
$ws.Cells.Item(2, 20).Value = 0.01699495368149812

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Sertad1"
$ws.Cells.Item(3, 3).Value = "Ar"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 22.02847633333333
$ws.Cells.Item(3, 8).Value = 66.085429
$ws.Cells.Item(3, 9).Value = 0.3378777065804683
$ws.Cells.Item(3, 10).Value = 0.3378777065804683
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 8.022254333333334
$ws.Cells.Item(3, 14).Value = 24.066763
$ws.Cells.Item(3, 15).Value = 0.5993413978297139
$ws.Cells.Item(3, 16).Value = 0.5993413978297139
$ws.Cells.Item(3, 17).Value = 176.7180397218141
$ws.Cells.Item(3, 18).Value = 1590.462357496327
$ws.Cells.Item(3, 19).Value = 0.2025040969574358
$ws.Cells.Item(3, 20).Value = 0.2025040969574358

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Sertad1"
$ws.Cells.Item(4, 3).Value = "Ar"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 22.02847633333333
$ws.Cells.Item(4, 8).Value = 66.085429
$ws.Cells.Item(4, 9).Value = 0.3378777065804683
$ws.Cells.Item(4, 10).Value = 0.3378777065804683
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.689602333333333
$ws.Cells.Item(4, 14).Value = 14.068807
$ws.Cells.Item(4, 15).Value = 0.3503594751473832
$ws.Cells.Item(4, 16).Value = 0.3503594751473832
$ws.Cells.Item(4, 17).Value = 103.3047940125781
$ws.Cells.Item(4, 18).Value = 929.743146113203
$ws.Cells.Item(4, 19).Value = 0.1183786559415344
$ws.Cells.Item(4, 20).Value = 0.1183786559415344

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Sertad1"
$ws.Cells.Item(5, 3).Value = "Ar"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 16.174732
$ws.Cells.Item(5, 8).Value = 48.524196
$ws.Cells.Item(5, 9).Value = 0.2480916641721602
$ws.Cells.Item(5, 10).Value = 0.2480916641721602
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.6732596666666667
$ws.Cells.Item(5, 14).Value = 2.019779
$ws.Cells.Item(5, 15).Value = 0.05029912702290298
$ws.Cells.Item(5, 16).Value = 0.05029912702290298
$ws.Cells.Item(5, 17).Value = 10.88979467474267
$ws.Cells.Item(5, 18).Value = 98.00815207268401
$ws.Cells.Item(5, 19).Value = 0.01247879412951887
$ws.Cells.Item(5, 20).Value = 0.01247879412951887

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Sertad1"
$ws.Cells.Item(6, 3).Value = "Ar"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 16.174732
$ws.Cells.Item(6, 8).Value = 48.524196
$ws.Cells.Item(6, 9).Value = 0.2480916641721602
$ws.Cells.Item(6, 10).Value = 0.2480916641721602
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 8.022254333333334
$ws.Cells.Item(6, 14).Value = 24.066763
$ws.Cells.Item(6, 15).Value = 0.5993413978297139
$ws.Cells.Item(6, 16).Value = 0.5993413978297139
$ws.Cells.Item(6, 17).Value = 129.7578138775054
$ws.Cells.Item(6, 18).Value = 1167.820324897548
$ws.Cells.Item(6, 19).Value = 0.1486916047948424
$ws.Cells.Item(6, 20).Value = 0.1486916047948424

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Sertad1"
$ws.Cells.Item(7, 3).Value = "Ar"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 16.174732
$ws.Cells.Item(7, 8).Value = 48.524196
$ws.Cells.Item(7, 9).Value = 0.2480916641721602
$ws.Cells.Item(7, 10).Value = 0.2480916641721602
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.689602333333333
$ws.Cells.Item(7, 14).Value = 14.068807
$ws.Cells.Item(7, 15).Value = 0.3503594751473832
$ws.Cells.Item(7, 16).Value = 0.3503594751473832
$ws.Cells.Item(7, 17).Value = 75.85306092824135
$ws.Cells.Item(7, 18).Value = 682.677548354172
$ws.Cells.Item(7, 19).Value = 0.08692126524779889
$ws.Cells.Item(7, 20).Value = 0.08692126524779889

$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Sertad1"
$ws.Cells.Item(8, 3).Value = "Ar"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 16.59113633333333
$ws.Cells.Item(8, 8).Value = 49.773409
$ws.Cells.Item(8, 9).Value = 0.254478567153005
$ws.Cells.Item(8, 10).Value = 0.2544785671530049
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.6732596666666667
$ws.Cells.Item(8, 14).Value = 2.019779
$ws.Cells.Item(8, 15).Value = 0.05029912702290298
$ws.Cells.Item(8, 16).Value = 0.05029912702290298
$ws.Cells.Item(8, 17).Value = 11.17014291740122
$ws.Cells.Item(8, 18).Value = 100.531286256611
$ws.Cells.Item(8, 19).Value = 0.01280004977383534
$ws.Cells.Item(8, 20).Value = 0.01280004977383534

$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Sertad1"
$ws.Cells.Item(9, 3).Value = "Ar"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 16.59113633333333
$ws.Cells.Item(9, 8).Value = 49.773409
$ws.Cells.Item(9, 9).Value = 0.254478567153005
$ws.Cells.Item(9, 10).Value = 0.2544785671530049
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 8.022254333333334
$ws.Cells.Item(9, 14).Value = 24.066763
$ws.Cells.Item(9, 15).Value = 0.5993413978297139
$ws.Cells.Item(9, 16).Value = 0.5993413978297139
$ws.Cells.Item(9, 17).Value = 133.0983153450075
$ws.Cells.Item(9, 18).Value = 1197.884838105067
$ws.Cells.Item(9, 19).Value = 0.1525195401551847
$ws.Cells.Item(9, 20).Value = 0.1525195401551847

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Sertad1"
$ws.Cells.Item(10, 3).Value = "Ar"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 16.59113633333333
$ws.Cells.Item(10, 8).Value = 49.773409
$ws.Cells.Item(10, 9).Value = 0.254478567153005
$ws.Cells.Item(10, 10).Value = 0.2544785671530049
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.689602333333333
$ws.Cells.Item(10, 14).Value = 14.068807
$ws.Cells.Item(10, 15).Value = 0.3503594751473832
$ws.Cells.Item(10, 16).Value = 0.3503594751473832
$ws.Cells.Item(10, 17).Value = 77.80583166145145
$ws.Cells.Item(10, 18).Value = 700.252484953063
$ws.Cells.Item(10, 19).Value = 0.08915897722398494
$ws.Cells.Item(10, 20).Value = 0.08915897722398491

$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Sertad1"
$ws.Cells.Item(11, 3).Value = "Ar"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 10.40225133333333
$ws.Cells.Item(11, 8).Value = 31.206754
$ws.Cells.Item(11, 9).Value = 0.1595520620943666
$ws.Cells.Item(11, 10).Value = 0.1595520620943666
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.6732596666666667
$ws.Cells.Item(11, 14).Value = 2.019779
$ws.Cells.Item(11, 15).Value = 0.05029912702290298
$ws.Cells.Item(11, 16).Value = 0.05029912702290298
$ws.Cells.Item(11, 17).Value = 7.003416265262889
$ws.Cells.Item(11, 18).Value = 63.030746387366
$ws.Cells.Item(11, 19).Value = 0.008025329438050649
$ws.Cells.Item(11, 20).Value = 0.008025329438050649

$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Sertad1"
$ws.Cells.Item(12, 3).Value = "Ar"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 10.40225133333333
$ws.Cells.Item(12, 8).Value = 31.206754
$ws.Cells.Item(12, 9).Value = 0.1595520620943666
$ws.Cells.Item(12, 10).Value = 0.1595520620943666
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 8.022254333333334
$ws.Cells.Item(12, 14).Value = 24.066763
$ws.Cells.Item(12, 15).Value = 0.5993413978297139
$ws.Cells.Item(12, 16).Value = 0.5993413978297139
$ws.Cells.Item(12, 17).Value = 83.44950583525578
$ws.Cells.Item(12, 18).Value = 751.045552517302
$ws.Cells.Item(12, 19).Value = 0.09562615592225097
$ws.Cells.Item(12, 20).Value = 0.09562615592225097

$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Sertad1"
$ws.Cells.Item(13, 3).Value = "Ar"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 10.40225133333333
$ws.Cells.Item(13, 8).Value = 31.206754
$ws.Cells.Item(13, 9).Value = 0.1595520620943666
$ws.Cells.Item(13, 10).Value = 0.1595520620943666
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 4.689602333333333
$ws.Cells.Item(13, 14).Value = 14.068807
$ws.Cells.Item(13, 15).Value = 0.3503594751473832
$ws.Cells.Item(13, 16).Value = 0.3503594751473832
$ws.Cells.Item(13, 17).Value = 48.78242212471977
$ws.Cells.Item(13, 18).Value = 439.0417991224779
$ws.Cells.Item(13, 19).Value = 0.05590057673406497
$ws.Cells.Item(13, 20).Value = 0.05590057673406497
